# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Integral" colour scheme (used by the slide master,
#                            i.e. every slide in the deck)
#   ppt/theme/theme2.xml -> "Office Theme" colour scheme (used only by the
#                            notes master)
#
# The authored change swaps the two themes' contents: the slide-facing theme
# (theme1.xml) becomes the standard "Office" palette, and the notes-master
# theme becomes the old "Integral" palette. The font scheme and format
# (fill/line/effect) scheme are identical between the two themes already, so
# the only substantive, visible difference is the 12 colour-scheme entries
# used throughout the slides.
#
# Apply the new "Office" theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) to the presentation's shared theme colour scheme so every slide
# picks up the standard Office colours instead of the Integral ones.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

$officeColors = @(
    0,          # 1  dk1      000000
    16777215,   # 2  lt1      FFFFFF
    6968388,    # 3  dk2      44546A
    15132391,   # 4  lt2      E7E6E6
    13998939,   # 5  accent1  5B9BD5
    3243501,    # 6  accent2  ED7D31
    10855845,   # 7  accent3  A5A5A5
    49407,      # 8  accent4  FFC000
    12874308,   # 9  accent5  4472C4
    4697456,    # 10 accent6  70AD47
    12673797,   # 11 hlink    0563C1
    7491477     # 12 folHlink 954F72
)

for ($i = 1; $i -le $cs.Count; $i++) {
    $cs.Item($i).RGB = $officeColors[$i - 1]
}
